# Make indst/IELC data-driven (#8)
#
# Inserts a new row into the "Key to Variables" sheet, right after the
# existing "indst" / FYIESM row (row 167), for the new "IELC" acronym
# (Industrial Equipment Logit Coefficient). All rows from the old row 168
# onward shift down by one, which Excel's Rows.Insert() handles natively
# (and also carries the formatting down from the row above, matching the
# target styling for the new row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")
$ws.Activate()

# Insert a new blank row at row 168 (pushes the old 168..277 down to 169..278).
$ws.Rows.Item(168).Insert()

# Populate the new row: Top Level Folder / Acronym / Meaning / Importance column.
$ws.Range("A168").Value = "indst"
$ws.Range("B168").Value = "IELC"
$ws.Range("C168").Value = "Industrial Equipment Logit Coefficient"
$ws.Range("F168").Value = "to be determined via calibration"

# Match the author's final selection/cursor position.
$ws.Range("A168").Select()
